$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C","D","E","F","J","K","L","M","N")

$data = @{
    2 = @(1.007517189882854, 1.010163023887674, 1.010173097929085, 1.00570656476623, 1.012788370573424, 1.013034046953649, 1.013044090525428, 1.008591131989706, 1.008322484023764)
    3 = @(1.008777202253876, 1.011333579330226, 1.011251613067969, 1.007621938967664, 1.013678231869978, 1.014007589187392, 1.013925852672849, 1.010306393046803, 1.008628179638942)
    4 = @(1.009591924433954, 1.012090770948063, 1.011949264209807, 1.008860251563744, 1.014253038152566, 1.014636742348009, 1.014495612020765, 1.011414844666994, 1.008825338359081)
    5 = @(1.009934297215638, 1.012409041674543, 1.012242507248352, 1.009380597362823, 1.014494452321081, 1.014901051889394, 1.014734950459729, 1.011880505766299, 1.008908069990505)
    6 = @(1.009991775252639, 1.01246247775598, 1.012291741235622, 1.009467951983472, 1.014534973177194, 1.014945419794747, 1.01477512547815, 1.011958673196527, 1.00892195199838)
    7 = @(1.009596499765874, 1.012095023899101, 1.011953182728163, 1.008867205379285, 1.014256264857257, 1.014640274793233, 1.014498810809175, 1.01142106814292, 1.008826444426729)
    8 = @(1.007943142064123, 1.010558669789744, 1.010537634729583, 1.006354102649743, 1.013089311152642, 1.013363226625048, 1.01334225360938, 1.009171117304428, 1.008425929781417)
    9 = @(1.005024977690948, 1.007849435501135, 1.008041425703835, 1.001916947197382, 1.011025225612805, 1.011106664175804, 1.011297989151957, 1.005194869748531, 1.007715175401642)
    10 = @(1.003076043972148, 1.006041700691043, 1.006375848360975, 0.9989521099072129, 1.009643729766757, 1.009597866282021, 1.009930740491634, 1.002535517254903, 1.007237918939261)
    11 = @(1.002231237537961, 1.005258501728548, 1.005654247396488, 0.9976665156659695, 1.00904418756963, 1.008943439922909, 1.009337617796227, 1.001381798509036, 1.007030436511824)
    12 = @(1.001917296819954, 1.004967516128488, 1.005386149240515, 0.9971887018993029, 1.008821284461942, 1.008700186052415, 1.009117137245418, 1.000952910667624, 1.006953242781246)
    13 = @(1.001984644691727, 1.005029936802071, 1.00544366014131, 0.9972912078179816, 1.008869107346712, 1.008752372633971, 1.00916443875357, 1.001044924427577, 1.006969806798484)
    14 = @(1.002205290031676, 1.005234450237169, 1.005632087636968, 0.9976270253773062, 1.009025766565622, 1.008923335992049, 1.009319396262325, 1.001346353635604, 1.007024058221722)
    15 = @(1.002341217912531, 1.005360448208953, 1.005748175508191, 0.997833895159757, 1.009122261995732, 1.00902864939447, 1.009414848270601, 1.001532027921628, 1.007057467644408)
    16 = @(1.003132091409819, 1.006093669388492, 1.00642372997058, 0.9990373912417067, 1.00968349069489, 1.009641274701102, 1.009970080692238, 1.002612038143036, 1.007251671330573)
    17 = @(1.003627938879378, 1.006553479023196, 1.006847378814024, 0.9997918190738038, 1.010035171399794, 1.010025258865606, 1.010318067441776, 1.003288900350847, 1.007373267829813)
    18 = @(1.003917071380771, 1.00682163628692, 1.007094448185242, 1.000231691490648, 1.010240171289094, 1.010249123836221, 1.010520936836497, 1.003683491195755, 1.007444113279416)
    19 = @(1.004015643525776, 1.006913064020034, 1.007178686161864, 1.000381647839342, 1.010310049115234, 1.010325438077179, 1.010590092297764, 1.003818001053591, 1.007468256260234)
    20 = @(1.003574748174463, 1.006504150171865, 1.006801929245292, 0.9997108940757595, 1.009997452816173, 1.009984072046661, 1.010280742685765, 1.003216301416195, 1.007360229932181)
    21 = @(1.002140319446498, 1.005174228116868, 1.005576602219858, 0.9975281435260553, 1.008979640063874, 1.00887299629484, 1.009273769817297, 1.001257599897155, 1.007008086004692)
    22 = @(1.001237612647247, 1.004337642242208, 1.004805821458389, 0.9961540939351059, 1.008338504223621, 1.008173428531071, 1.008639669918193, 1.000024081236138, 1.006785952305349)
    23 = @(1.001716234765744, 1.004781172838824, 1.005214463276155, 0.9968826672004004, 1.008678497332354, 1.008544378016595, 1.00897591209132, 1.000678187788766, 1.006903778894003)
    24 = @(1.003598783033375, 1.00652643988628, 1.006822466066263, 0.9997474611321411, 1.010014496620956, 1.010002682926595, 1.010297608462273, 1.003249106403021, 1.007366121443386)
    25 = @(1.005779985857746, 1.008550097602331, 1.008686993834474, 1.00306518222948, 1.011559783695783, 1.011690802469281, 1.011827242669469, 1.006224271044895, 1.007899520777064)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

$wb.Save()
